$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# The "2 Player?" column (F) for rows 5-12 should be set to "Yes"
# instead of "No", reflecting that 2-player logic now works for
# these 2P-titled minigames.
$ws.Range("F5:F12").Value = "Yes"
